$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.800.31"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").Value = "2.634.27"
$ws.Range("E3").Value = "  +1.81%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.71"
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.37"
$ws.Range("E6").Value = "  +3.74%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.568"
$ws.Range("E8").Value = "  +0.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.57"
$ws.Range("E9").Value = "  +2.31%  "
$ws.Range("E10").Value = "  +2.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.337"
$ws.Range("E11").Value = "  +1.73%  "
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("D13").Value = "3.102.04"
$ws.Range("E13").Value = "  +1.92%  "
$ws.Range("D14").Value = "59.711.53"
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.97"
$ws.Range("E15").Value = "  +2.48%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000135"
$ws.Range("E16").Value = "  +1.69%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.621.79"
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "342.52"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("E19").Value = "  +2.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.22"
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.39"
$ws.Range("E21").Value = "  -0.33%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.59"
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("E24").Value = "  +1.91%  "
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.25"
$ws.Range("E27").Value = "  +2.78%  "
$ws.Range("D28").Value = "0.0₃0753"
$ws.Range("E28").Value = "  +5.26%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  +4.04%  "
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.00"
$ws.Range("E32").Value = "  +1.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "151.00"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.02"
$ws.Range("E34").Value = "  +1.72%  "
$ws.Range("E35").Value = "  +2.34%  "
$ws.Range("E36").Value = "  -0.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.839"
$ws.Range("E37").Value = "  +2.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.825"
$ws.Range("E38").Value = "  +2.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "292.38"
$ws.Range("E39").Value = "  +9.31%  "
$ws.Range("E40").Value = "  +2.02%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("E42").Value = "  +0.99%  "
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0533"
$ws.Range("E45").Value = "  +3.78%  "
$ws.Range("D46").Value = "1.970.49"
$ws.Range("E46").Value = "  +0.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0225"
$ws.Range("E47").Value = "  +2.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.47"
$ws.Range("E48").Value = "  +2.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.56"
$ws.Range("E49").Value = "  +2.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.12"
$ws.Range("E50").Value = "  +0.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.74"
$ws.Range("E51").Value = "  -0.12%  "
